# Update "想去人数" (F) and "最低票价" (G) figures across the four sheets
# to match the newly scraped data (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1781
$ws.Range("G3").Value = 110
$ws.Range("F5").Value = 629
$ws.Range("F6").Value = 1156
$ws.Range("F7").Value = 1588
$ws.Range("F8").Value = 176
$ws.Range("F9").Value = 176
$ws.Range("F11").Value = 1525
$ws.Range("F13").Value = 685
$ws.Range("F14").Value = 1844
$ws.Range("F15").Value = 1832
$ws.Range("F16").Value = 892
$ws.Range("F17").Value = 304
$ws.Range("F19").Value = 1513
$ws.Range("F20").Value = 309
$ws.Range("F22").Value = 30
$ws.Range("F23").Value = 1302
$ws.Range("F24").Value = 424
$ws.Range("F25").Value = 518
$ws.Range("F26").Value = 195
$ws.Range("F27").Value = 6840
$ws.Range("F28").Value = 5444
$ws.Range("F29").Value = 777
$ws.Range("F31").Value = 1717
$ws.Range("F33").Value = 245

# --- Sheet: 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 26

# --- Sheet: 本地生活 (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 51

# --- Sheet: 全部类型 (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 51
$ws.Range("F4").Value = 26
$ws.Range("F5").Value = 1781
$ws.Range("G5").Value = 110
$ws.Range("F7").Value = 629
$ws.Range("F8").Value = 1156
$ws.Range("F9").Value = 1588
$ws.Range("F10").Value = 176
$ws.Range("F11").Value = 176
$ws.Range("F14").Value = 1525
$ws.Range("F16").Value = 685
$ws.Range("F17").Value = 1844
$ws.Range("F18").Value = 1832
$ws.Range("F19").Value = 892
$ws.Range("F20").Value = 304
$ws.Range("F22").Value = 1513
$ws.Range("F23").Value = 309
$ws.Range("F26").Value = 30
$ws.Range("F28").Value = 1302
$ws.Range("F29").Value = 424
$ws.Range("F30").Value = 518
$ws.Range("F31").Value = 195
$ws.Range("F32").Value = 6840
$ws.Range("F33").Value = 5444
$ws.Range("F34").Value = 777
$ws.Range("F36").Value = 1717
$ws.Range("F40").Value = 245

$wb.Save()
